$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("K50").Value = "Myoxocephalus quadricornis"
$ws.Range("A50").Style = "Good"
$ws.Range("B50").Style = "Good"
$ws.Range("K50").Style = "Good"
